$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 880
$ws.Range("J40").Value = 796.6667
$ws.Range("L40").Value = 796.6667
$ws.Range("N40").Value = -1146.6667
$ws.Range("H76").Value = 3321.9375
$ws.Range("I76").Value = 3076.25
$ws.Range("J76").Value = 4059
$ws.Range("K76").Value = 3076.25
$ws.Range("L76").Value = 4059
$ws.Range("M76").Value = -2761.25
$ws.Range("N76").Value = -4689
$ws.Range("H79").Value = 3321.9375
$ws.Range("I79").Value = 3076.25
$ws.Range("J79").Value = 4059
$ws.Range("K79").Value = 3076.25
$ws.Range("L79").Value = 4059
$ws.Range("M79").Value = -1984.25
$ws.Range("N79").Value = -6243
$ws.Range("H97").Value = 10000
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -30992
$ws.Range("H98").Value = 1923.6875
$ws.Range("I98").Value = 2198.2222
$ws.Range("K98").Value = 2198.2222
$ws.Range("M98").Value = -700.2222000000002
$ws.Range("H113").Value = 1781.6666
$ws.Range("I113").Value = 1675.5555
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1675.5555
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1578.4445
$ws.Range("N113").Value = -8608
$ws.Range("H118").Value = 1232.7222
$ws.Range("I118").Value = 1356
$ws.Range("J118").Value = 1185.3077
$ws.Range("K118").Value = 4068
$ws.Range("L118").Value = 3555.9231
$ws.Range("M118").Value = -2411
$ws.Range("N118").Value = -6869.9231
$ws.Range("H122").Value = 1923.6875
$ws.Range("I122").Value = 2198.2222
$ws.Range("K122").Value = 6594.6666
$ws.Range("M122").Value = -4144.6666
$ws.Range("H137").Value = 9477.75
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 12303.667
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 36911.001
$ws.Range("M137").Value = -450
$ws.Range("N137").Value = -42011.001
$ws.Range("H138").Value = 1579
$ws.Range("I138").Value = 679.0714
$ws.Range("J138").Value = 2103.9583
$ws.Range("K138").Value = 2037.2142
$ws.Range("L138").Value = 6311.874899999999
$ws.Range("M138").Value = 3102.7858
$ws.Range("N138").Value = -16591.8749

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 288592.56
$ws.Range("I74").Value = 288592.56
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 288592.56
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -287718.56
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 288592.56
$ws.Range("I77").Value = 288592.56
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 1442962.8
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1438594.8
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 1523.7778
$ws.Range("I122").Value = 1281
$ws.Range("K122").Value = 3843
$ws.Range("M122").Value = -1393

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 27495.8
$ws.Range("J81").Value = 27495.8
$ws.Range("L81").Value = 27495.8
$ws.Range("N81").Value = -29617.8
$ws.Range("H84").Value = 27495.8
$ws.Range("J84").Value = 27495.8
$ws.Range("L84").Value = 82487.39999999999
$ws.Range("N84").Value = -93095.39999999999
$ws.Range("H99").Value = 1393.6364
$ws.Range("I99").Value = 1302.9259
$ws.Range("J99").Value = 1801.8334
$ws.Range("K99").Value = 1302.9259
$ws.Range("L99").Value = 1801.8334
$ws.Range("M99").Value = 195.0741
$ws.Range("N99").Value = -4797.8334
$ws.Range("H134").Value = 27836846
$ws.Range("I134").Value = 41668564
$ws.Range("J134").Value = 173412.33
$ws.Range("K134").Value = 125005692
$ws.Range("L134").Value = 520236.99
$ws.Range("M134").Value = -125003157
$ws.Range("N134").Value = -525306.99

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37391.566
$ws.Range("I31").Value = 58093.332
$ws.Range("J31").Value = 10220.5
$ws.Range("K31").Value = 58093.332
$ws.Range("L31").Value = 10220.5
$ws.Range("M31").Value = -57798.332
$ws.Range("N31").Value = -10810.5
$ws.Range("H34").Value = 37391.566
$ws.Range("I34").Value = 58093.332
$ws.Range("J34").Value = 10220.5
$ws.Range("K34").Value = 58093.332
$ws.Range("L34").Value = 10220.5
$ws.Range("M34").Value = -57891.332
$ws.Range("N34").Value = -10624.5
$ws.Range("H122").Value = 1127.4286
$ws.Range("I122").Value = 1098.6666
$ws.Range("K122").Value = 3295.9998
$ws.Range("M122").Value = -845.9998000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 28893136
$ws.Range("J55").Value = 36115996
$ws.Range("L55").Value = 108347988
$ws.Range("N55").Value = -108348342
$ws.Range("H117").Value = 943.4545000000001
$ws.Range("I117").Value = 416.66666
$ws.Range("J117").Value = 1141
$ws.Range("K117").Value = 1249.99998
$ws.Range("L117").Value = 3423
$ws.Range("M117").Value = 2192.00002
$ws.Range("N117").Value = -10307
$ws.Range("H129").Value = 1491
$ws.Range("J129").Value = 1625.65
$ws.Range("L129").Value = 4876.950000000001
$ws.Range("N129").Value = -14876.95

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 14805.634
$ws.Range("I102").Value = 5360.44
$ws.Range("K102").Value = 5360.44
$ws.Range("M102").Value = -3738.44
$ws.Range("H122").Value = 1083.1666
$ws.Range("J122").Value = 999.8
$ws.Range("L122").Value = 2999.4
$ws.Range("N122").Value = -7899.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2512.6553
$ws.Range("I7").Value = 2967.75
$ws.Range("K7").Value = 2967.75
$ws.Range("M7").Value = -2855.75
$ws.Range("H40").Value = 75057.07000000001
$ws.Range("I40").Value = 2696.2856
$ws.Range("J40").Value = 147417.86
$ws.Range("K40").Value = 2696.2856
$ws.Range("L40").Value = 147417.86
$ws.Range("M40").Value = -2560.2856
$ws.Range("N40").Value = -147689.86
$ws.Range("H122").Value = 2574.9644
$ws.Range("I122").Value = 2400.2354
$ws.Range("J122").Value = 2845
$ws.Range("K122").Value = 7200.706200000001
$ws.Range("L122").Value = 8535
$ws.Range("M122").Value = -4750.706200000001
$ws.Range("N122").Value = -13435
$ws.Range("H126").Value = 2512.6553
$ws.Range("I126").Value = 2967.75
$ws.Range("K126").Value = 8903.25
$ws.Range("M126").Value = -6433.25
$ws.Range("H136").Value = 1671.6562
$ws.Range("I136").Value = 926.7059
$ws.Range("J136").Value = 2515.9333
$ws.Range("K136").Value = 2780.1177
$ws.Range("L136").Value = 7547.7999
$ws.Range("M136").Value = -230.1177000000002
$ws.Range("N136").Value = -12647.7999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5020.8945
$ws.Range("I122").Value = 3889.7
$ws.Range("J122").Value = 6277.778
$ws.Range("K122").Value = 11669.1
$ws.Range("L122").Value = 18833.334
$ws.Range("M122").Value = -9219.099999999999
$ws.Range("N122").Value = -23733.334
